$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H132").Value = 2098.5
$ws.Range("I132").Value = 1110.875
$ws.Range("K132").Value = 3332.625
$ws.Range("M132").Value = -802.625
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2666.25
$ws.Range("I2").Value = 2666.25
$ws.Range("K2").Value = 2666.25
$ws.Range("M2").Value = -2553.25
$ws.Range("H32").Value = 8003
$ws.Range("I32").Value = 8003
$ws.Range("K32").Value = 8003
$ws.Range("M32").Value = -7716
$ws.Range("H45").Value = 2570.6667
$ws.Range("I45").Value = 2570.6667
$ws.Range("K45").Value = 2570.6667
$ws.Range("M45").Value = -2193.6667
$ws.Range("H102").Value = 2874.4
$ws.Range("I102").Value = 2874.4
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2874.4
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1252.4
$ws.Range("N102").ClearContents()
$ws.Range("H116").Value = 2666.25
$ws.Range("I116").Value = 2666.25
$ws.Range("K116").Value = 2666.25
$ws.Range("M116").Value = -372.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2666.25
$ws.Range("I3").Value = 2666.25
$ws.Range("K3").Value = 2666.25
$ws.Range("M3").Value = -2552.25
$ws.Range("H22").Value = 470
$ws.Range("I22").Value = 470
$ws.Range("K22").Value = 470
$ws.Range("M22").Value = -297
$ws.Range("H105").Value = 10298.6
$ws.Range("I105").Value = 10298.6
$ws.Range("K105").Value = 10298.6
$ws.Range("M105").Value = -8551.6
$ws.Range("H134").Value = 5087.143
$ws.Range("I134").Value = 5568.3335
$ws.Range("J134").Value = 2200
$ws.Range("K134").Value = 16705.0005
$ws.Range("L134").Value = 6600
$ws.Range("M134").Value = -14170.0005
$ws.Range("N134").Value = -11670
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 78795
$ws.Range("J68").Value = 78795
$ws.Range("L68").Value = 78795
$ws.Range("N68").Value = -80293
$ws.Range("H71").Value = 78795
$ws.Range("J71").Value = 78795
$ws.Range("L71").Value = 236385
$ws.Range("N71").Value = -243873
$ws.Range("H74").Value = 48996.285
$ws.Range("J74").Value = 48996.285
$ws.Range("L74").Value = 48996.285
$ws.Range("N74").Value = -50744.285
$ws.Range("H77").Value = 48996.285
$ws.Range("J77").Value = 48996.285
$ws.Range("L77").Value = 146988.855
$ws.Range("N77").Value = -155724.855
$ws.Range("H122").Value = 2343.9
$ws.Range("I122").Value = 2088
$ws.Range("J122").Value = 2599.8
$ws.Range("K122").Value = 6264
$ws.Range("L122").Value = 7799.400000000001
$ws.Range("M122").Value = -3814
$ws.Range("N122").Value = -12699.4
$ws.Range("H132").Value = 104095.6
$ws.Range("I132").Value = 145708.28
$ws.Range("J132").Value = 6999.3335
$ws.Range("K132").Value = 437124.84
$ws.Range("L132").Value = 20998.0005
$ws.Range("M132").Value = -434594.84
$ws.Range("N132").Value = -26058.0005
$ws.Range("H134").Value = 8600
$ws.Range("I134").Value = 9750
$ws.Range("K134").Value = 29250
$ws.Range("M134").Value = -26715
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 4927.875
$ws.Range("I81").Value = 3974.3333
$ws.Range("J81").Value = 5500
$ws.Range("K81").Value = 11922.9999
$ws.Range("L81").Value = 16500
$ws.Range("M81").Value = -10799.9999
$ws.Range("N81").Value = -18746
$ws.Range("H84").Value = 4927.875
$ws.Range("I84").Value = 3974.3333
$ws.Range("J84").Value = 5500
$ws.Range("K84").Value = 35768.9997
$ws.Range("L84").Value = 49500
$ws.Range("M84").Value = -30152.9997
$ws.Range("N84").Value = -60732
$ws.Range("H137").Value = 7499.75
$ws.Range("J137").Value = 4999.6665
$ws.Range("L137").Value = 14998.9995
$ws.Range("N137").Value = -25198.9995
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8333
$ws.Range("I70").Value = 8333
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 8333
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -8063
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 8333
$ws.Range("I73").Value = 8333
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 8333
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -7397
$ws.Range("N73").ClearContents()
$ws.Range("H97").Value = 2115.125
$ws.Range("I97").Value = 2345.8572
$ws.Range("K97").Value = 2345.8572
$ws.Range("M97").Value = -1849.8572
$ws.Range("H113").Value = 3975
$ws.Range("I113").Value = 3975
$ws.Range("K113").Value = 3975
$ws.Range("M113").Value = -1805
$ws.Range("H122").Value = 2348.3
$ws.Range("I122").Value = 2260.375
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 6781.125
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -4331.125
$ws.Range("N122").Value = -13000
$ws.Range("H132").Value = 3199.4443
$ws.Range("J132").Value = 4999
$ws.Range("L132").Value = 14997
$ws.Range("N132").Value = -20057
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4273.727
$ws.Range("I82").Value = 2470.25
$ws.Range("K82").Value = 2470.25
$ws.Range("M82").Value = -2109.25
$ws.Range("H85").Value = 4273.727
$ws.Range("I85").Value = 2470.25
$ws.Range("K85").Value = 2470.25
$ws.Range("M85").Value = -1222.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 29000
$ws.Range("I64").Value = 29000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 29000
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -28752
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 29000
$ws.Range("I67").Value = 29000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 29000
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -28142
$ws.Range("N67").ClearContents()
$ws.Range("H107").Value = 673.3077
$ws.Range("I107").Value = 636.8
$ws.Range("K107").Value = 1910.4
$ws.Range("M107").Value = 9.600000000000136
$ws.Range("H122").Value = 3083.1667
$ws.Range("I122").Value = 3083.1667
$ws.Range("K122").Value = 9249.500100000001
$ws.Range("M122").Value = -6799.500100000001
